# Update buffer address sheet with info from SEEED mail.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 12: clarify the calibration-factor-1 coefficient note.
$ws.Range("I12").Value = "calculation coeff for calibration factor 1 (mail Seeed)"

# Row 14: add the calibration-factor-2 coefficient note (previously empty).
$ws.Range("I14").Value = "calculated coefficient for calibration factor 2 (mail Seeed)"

# Row 16: device identification register.
$ws.Range("I16").Value = "device identification ( mail seeed)"

# Row 17: device address register.
$ws.Range("I17").Value = "device address (mail Seeed)"

# Row 19: year of factory shipment.
$ws.Range("I19").Value = "year of factory shipment (mail Seeed)"

# Row 20: date of factory shipment.
$ws.Range("I20").Value = "date of factory shipment (mail Seeed) Day month ?"

# Row 21: identification label for the reset-to-factory-settings register.
$ws.Range("A21").Value = "resst to factory settings"

$ws.Range("I20").Select()
